$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: D10 date changes from 41671 to 41615 ---
$ws.Range("D10").Value = 41615

# --- Row 11: situation -> Desenvolvido, date instead of N/A ---
$ws.Range("C11").Value = "Desenvolvido"
$ws.Range("D11").Value = 41615

# --- Row 12: situation -> Desenvolvido, date instead of N/A ---
$ws.Range("C12").Value = "Desenvolvido"
$ws.Range("D12").Value = 41615

# --- Row 13: situation -> Desenvolvido, date instead of N/A ---
$ws.Range("C13").Value = "Desenvolvido"
$ws.Range("D13").Value = 41615

# --- Row 14 (new): defect about "procurador" ---
$ws.Range("A14").Value = "Problema ao selecionar um procurador no cadastro de patentes"
$ws.Range("B14").Value = "Defeito"
$ws.Range("C14").Value = "Desenvolvido"
$ws.Range("D14").Value = 41615

# --- Row 15 (new): header-like row for "Relatórios de patente" ---
$ws.Range("A15").Value = "Relatórios de patente"
$ws.Range("B15").Value = "Novo item"
$ws.Range("C15").Value = "Em desenvolvimento"
$ws.Range("D15").Value = "N/A"
